$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.899.42"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "1.636.24"
$ws.Range("E3").Value = "  +1.13%  "
$ws.Range("E4").Value = "  +0.77%  "
$ws.Range("D5").Value = "'215.23"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.78%  "
$ws.Range("D8").Value = "'28.83"
$ws.Range("E8").Value = "  -1.49%  "
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'0.0901"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("E12").Value = "  +1.05%  "
$ws.Range("D13").Value = "1.594.71"
$ws.Range("E13").Value = "  -1.58%  "
$ws.Range("D14").Value = "'0.586"
$ws.Range("E14").Value = "  +3.67%  "
$ws.Range("E15").Value = "  +6.58%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "29.910.26"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("D18").Value = "'64.60"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "'240.58"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("D22").Value = "'9.92"
$ws.Range("E22").Value = "  +3.45%  "
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("E24").Value = "  +3.65%  "
$ws.Range("D25").Value = "'157.36"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "'0.110"
$ws.Range("E27").Value = "  -0.54%  "
$ws.Range("D28").Value = "'6.62"
$ws.Range("E28").Value = "  +0.81%  "
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("E32").Value = "  +1.29%  "
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "1.422.89"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +3.20%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("E37").Value = "  -3.40%  "
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").Value = "'76.26"
$ws.Range("E40").Value = "  +10.27%  "
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("E45").Value = "  +0.83%  "
$ws.Range("E46").Value = "  -1.67%  "
$ws.Range("E47").Value = "  +0.99%  "
$ws.Range("E48").Value = "  -1.80%  "
$ws.Range("D49").Value = "'48.83"
$ws.Range("E49").Value = "  -8.84%  "
$ws.Range("D50").Value = "'92.84"
$ws.Range("E50").Value = "  +5.29%  "
$ws.Range("D51").Value = "0.0₆0110"
$ws.Range("E51").Value = "  +6.78%  "
